$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date column C bumps by one day (46076 -> 46077) for all data rows.
$ws.Range("C2").Value = 46077
$ws.Range("C3").Value = 46077
$ws.Range("C4").Value = 46077
$ws.Range("C5").Value = 46077
$ws.Range("C6").Value = 46077
$ws.Range("C7").Value = 46077
$ws.Range("C8").Value = 46077

# Rows 6 and 7 swap their Beteckning / Datum / Area values.
$ws.Range("A6").Value = "A 21472-2022"
$ws.Range("B6").Value = 44706
$ws.Range("G6").Value = 1.5

$ws.Range("A7").Value = "A 21888-2022"
$ws.Range("B7").Value = 44709
$ws.Range("G7").Value = 2.4
